$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26479.744
$ws.Range("I2").Value = 30327.354
$ws.Range("J2").Value = 316
$ws.Range("K2").Value = 30327.354
$ws.Range("L2").Value = 316
$ws.Range("M2").Value = -30214.354
$ws.Range("N2").Value = -542

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20275.295
$ws.Range("I32").Value = 4182.5835
$ws.Range("J32").Value = 73917.664
$ws.Range("K32").Value = 4182.5835
$ws.Range("L32").Value = 73917.664
$ws.Range("M32").Value = -3895.5835
$ws.Range("N32").Value = -74491.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 956
$ws.Range("I45").Value = 956
$ws.Range("K45").Value = 956
$ws.Range("M45").Value = -579

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 4042.2222
$ws.Range("I97").Value = 4042.2222
$ws.Range("K97").Value = 4042.2222
$ws.Range("M97").Value = -3546.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 26479.744
$ws.Range("I116").Value = 30327.354
$ws.Range("J116").Value = 316
$ws.Range("K116").Value = 30327.354
$ws.Range("L116").Value = 316
$ws.Range("M116").Value = -28033.354
$ws.Range("N116").Value = -4904

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 10692.75
$ws.Range("I122").Value = 12231.3
$ws.Range("K122").Value = 36693.89999999999
$ws.Range("M122").Value = -34243.89999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3822.9333
$ws.Range("I132").Value = 3590.5
$ws.Range("J132").Value = 4088.5715
$ws.Range("K132").Value = 10771.5
$ws.Range("L132").Value = 12265.7145
$ws.Range("M132").Value = -8241.5
$ws.Range("N132").Value = -17325.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26479.744
$ws.Range("I3").Value = 30327.354
$ws.Range("J3").Value = 316
$ws.Range("K3").Value = 30327.354
$ws.Range("L3").Value = 316
$ws.Range("M3").Value = -30213.354
$ws.Range("N3").Value = -544

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1524.3334
$ws.Range("J86").Value = 1516.2858
$ws.Range("L86").Value = 1516.2858
$ws.Range("N86").Value = -3762.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1524.3334
$ws.Range("J89").Value = 1516.2858
$ws.Range("L89").Value = 7581.429
$ws.Range("N89").Value = -18813.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 33000
$ws.Range("J98").Value = 33000
$ws.Range("L98").Value = 33000
$ws.Range("N98").Value = -38990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 241123.81
$ws.Range("I105").Value = 2831.75
$ws.Range("K105").Value = 2831.75
$ws.Range("M105").Value = -1084.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5221.615
$ws.Range("I134").Value = 4400
$ws.Range("J134").Value = 5925.857
$ws.Range("K134").Value = 13200
$ws.Range("L134").Value = 17777.571
$ws.Range("M134").Value = -10665
$ws.Range("N134").Value = -22847.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4502.1187
$ws.Range("I31").Value = 1342.9445
$ws.Range("J31").Value = 9446.913
$ws.Range("K31").Value = 1342.9445
$ws.Range("L31").Value = 9446.913
$ws.Range("M31").Value = -1047.9445
$ws.Range("N31").Value = -10036.913

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4502.1187
$ws.Range("I34").Value = 1342.9445
$ws.Range("J34").Value = 9446.913
$ws.Range("K34").Value = 1342.9445
$ws.Range("L34").Value = 9446.913
$ws.Range("M34").Value = -1140.9445
$ws.Range("N34").Value = -9850.913

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3172.4285
$ws.Range("I99").Value = 2085.6667
$ws.Range("J99").Value = 3987.5
$ws.Range("K99").Value = 2085.6667
$ws.Range("L99").Value = 3987.5
$ws.Range("M99").Value = -587.6667000000002
$ws.Range("N99").Value = -6983.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3172.4285
$ws.Range("I126").Value = 2085.6667
$ws.Range("J126").Value = 3987.5
$ws.Range("K126").Value = 6257.000100000001
$ws.Range("L126").Value = 11962.5
$ws.Range("M126").Value = -3787.000100000001
$ws.Range("N126").Value = -16902.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 1000
$ws.Range("J93").Value = 1000
$ws.Range("L93").Value = 3000
$ws.Range("N93").Value = -6744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5241.4473
$ws.Range("I70").Value = 5316.3335
$ws.Range("K70").Value = 5316.3335
$ws.Range("M70").Value = -5046.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5241.4473
$ws.Range("I73").Value = 5316.3335
$ws.Range("K73").Value = 5316.3335
$ws.Range("M73").Value = -4380.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2612.7083
$ws.Range("I102").Value = 2468.7368
$ws.Range("J102").Value = 3159.8
$ws.Range("K102").Value = 2468.7368
$ws.Range("L102").Value = 3159.8
$ws.Range("M102").Value = -846.7368000000001
$ws.Range("N102").Value = -6403.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1645.3
$ws.Range("I122").Value = 1477.6666
$ws.Range("J122").Value = 1896.75
$ws.Range("K122").Value = 4432.9998
$ws.Range("L122").Value = 5690.25
$ws.Range("M122").Value = -1982.9998
$ws.Range("N122").Value = -10590.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 11290.667
$ws.Range("J123").Value = 11290.667
$ws.Range("L123").Value = 11290.667
$ws.Range("N123").Value = -16190.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2393.1162
$ws.Range("I126").Value = 1990.6666
$ws.Range("J126").Value = 2777.2727
$ws.Range("K126").Value = 5971.9998
$ws.Range("L126").Value = 8331.8181
$ws.Range("M126").Value = -3501.9998
$ws.Range("N126").Value = -13271.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3111.52
$ws.Range("I7").Value = 2999.8
$ws.Range("J7").Value = 3139.45
$ws.Range("K7").Value = 2999.8
$ws.Range("L7").Value = 3139.45
$ws.Range("M7").Value = -2887.8
$ws.Range("N7").Value = -3363.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3323.077
$ws.Range("I40").Value = 2783.3333
$ws.Range("J40").Value = 3485
$ws.Range("K40").Value = 2783.3333
$ws.Range("L40").Value = 3485
$ws.Range("M40").Value = -2647.3333
$ws.Range("N40").Value = -3757

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3153.2144
$ws.Range("I122").Value = 1542.8572
$ws.Range("J122").Value = 3690
$ws.Range("K122").Value = 4628.571599999999
$ws.Range("L122").Value = 11070
$ws.Range("M122").Value = -2178.571599999999
$ws.Range("N122").Value = -15970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3111.52
$ws.Range("I126").Value = 2999.8
$ws.Range("J126").Value = 3139.45
$ws.Range("K126").Value = 8999.400000000001
$ws.Range("L126").Value = 9418.349999999999
$ws.Range("M126").Value = -6529.400000000001
$ws.Range("N126").Value = -14358.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3283.3901
$ws.Range("I132").Value = 1896.8
$ws.Range("K132").Value = 5690.4
$ws.Range("M132").Value = -3160.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1613.7073
$ws.Range("I132").Value = 1108.0312
$ws.Range("K132").Value = 3324.0936
$ws.Range("M132").Value = -794.0935999999997

Write-Output "Applied all cell updates"